# Update the dSF (column F) values for the listed rows, per the
# "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -3
    5  = -2
    6  = -3
    7  = 2
    8  = -1
    9  = -4
    12 = -4
    13 = -1
    14 = 1
    15 = -1
    17 = 3
    18 = 1
    19 = 2
    20 = -6
    21 = -1
    22 = -3
    24 = -4
    25 = 3
    26 = -2
    28 = 2
    29 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
